$d = $word.ActiveDocument

# 1. Remove the duplicated bold "Play Free Bounding Luck Slot Game"
#    paragraph that sits just before the italic blurb near the end of the
#    document (do this first, while the text is still unique, before we
#    add a second "Meta description" paragraph containing similar text).
$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "Play Free Bounding Luck Slot Game`r") {
        $p.Range.Delete()
        break
    }
}

# 2. Replace the old italic blurb text with the new image-prompt text.
$old = "Read our review of Bounding Luck, an Asian-themed slot game with 1,024 paylines, high volatility, and bonus features available to play for free."
$new = "Prompt: Create a feature image for Bounding Luck that showcases the fun and exciting theme of the game. The image should be in a cartoon style and should prominently feature a happy Maya warrior wearing glasses. The main colors should be vibrant and eye-catching, preferably with an Asian-inspired color palette. The image should also include elements of the game, such as the bunny and the ethnic artifacts that can be found on the reels. Overall, the image should evoke a sense of adventure and luck, inviting players to join the Maya warrior on their journey to win big in Bounding Luck."
$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null

# 3. Insert a new "Meta description" paragraph right after the first
#    (Heading1) paragraph "Play Free Bounding Luck Slot Game".
$p1 = $d.Paragraphs.Item(1)
$p1.Range.InsertParagraphAfter() | Out-Null
$p2 = $d.Paragraphs.Item(2)

$metaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:r/>' +
           '<w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>' +
           '<w:r><w:t>: Read our review of Bounding Luck, an Asian-themed slot game with 1,024 paylines, high volatility, and bonus features available to play for free.</w:t></w:r>' +
           '</w:p>'
$p2.Range.InsertXML($metaXml)
